$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1514.3158
$ws.Range("I28").Value = 859.125
$ws.Range("J28").Value = 5008.6665
$ws.Range("K28").Value = 859.125
$ws.Range("L28").Value = 5008.6665
$ws.Range("M28").Value = -374.125
$ws.Range("N28").Value = -5978.6665
$ws.Range("H53").Value = 671.4737
$ws.Range("I53").Value = 1498.5
$ws.Range("J53").Value = 70
$ws.Range("K53").Value = 1498.5
$ws.Range("L53").Value = 70
$ws.Range("M53").Value = -861.5
$ws.Range("N53").Value = -1344
$ws.Range("H62").Value = 1827.3636
$ws.Range("I62").Value = 1429
$ws.Range("J62").Value = 2524.5
$ws.Range("K62").Value = 1429
$ws.Range("L62").Value = 2524.5
$ws.Range("M62").Value = -805
$ws.Range("N62").Value = -3772.5
$ws.Range("H65").Value = 1827.3636
$ws.Range("I65").Value = 1429
$ws.Range("J65").Value = 2524.5
$ws.Range("K65").Value = 7145
$ws.Range("L65").Value = 12622.5
$ws.Range("M65").Value = -4025
$ws.Range("N65").Value = -18862.5
$ws.Range("H92").Value = 1229.1052
$ws.Range("I92").Value = 1356.0588
$ws.Range("J92").Value = 150
$ws.Range("K92").Value = 1356.0588
$ws.Range("L92").Value = 150
$ws.Range("M92").Value = -108.0588
$ws.Range("N92").Value = -2646
$ws.Range("H96").Value = 6829.6665
$ws.Range("I96").Value = 7990.75
$ws.Range("J96").Value = 4507.5
$ws.Range("K96").Value = 23972.25
$ws.Range("L96").Value = 13522.5
$ws.Range("M96").Value = -22599.25
$ws.Range("N96").Value = -16268.5
$ws.Range("H98").Value = 741.0909
$ws.Range("I98").Value = 628.2778
$ws.Range("J98").Value = 1248.75
$ws.Range("K98").Value = 628.2778
$ws.Range("L98").Value = 1248.75
$ws.Range("M98").Value = 869.7222
$ws.Range("N98").Value = -4244.75
$ws.Range("H100").Value = 2251.1428
$ws.Range("I100").Value = 1962.5
$ws.Range("J100").Value = 2636
$ws.Range("K100").Value = 1962.5
$ws.Range("L100").Value = 2636
$ws.Range("M100").Value = -1421.5
$ws.Range("N100").Value = -3718
$ws.Range("H113").Value = 2001.2821
$ws.Range("I113").Value = 1902.5714
$ws.Range("J113").Value = 2116.4443
$ws.Range("K113").Value = 1902.5714
$ws.Range("L113").Value = 2116.4443
$ws.Range("M113").Value = 1351.4286
$ws.Range("N113").Value = -8624.444299999999
$ws.Range("H116").Value = 2725.9167
$ws.Range("I116").Value = 2221
$ws.Range("J116").Value = 3086.5715
$ws.Range("K116").Value = 2221
$ws.Range("L116").Value = 3086.5715
$ws.Range("M116").Value = 1221
$ws.Range("N116").Value = -9970.5715
$ws.Range("H122").Value = 741.0909
$ws.Range("I122").Value = 628.2778
$ws.Range("J122").Value = 1248.75
$ws.Range("K122").Value = 1884.8334
$ws.Range("L122").Value = 3746.25
$ws.Range("M122").Value = 565.1666
$ws.Range("N122").Value = -8646.25
$ws.Range("H132").Value = 2594.7812
$ws.Range("I132").Value = 2011.6786
$ws.Range("J132").Value = 6676.5
$ws.Range("K132").Value = 6035.0358
$ws.Range("L132").Value = 20029.5
$ws.Range("M132").Value = -3505.0358
$ws.Range("N132").Value = -25089.5
$ws.Range("H135").Value = 1445.4615
$ws.Range("I135").Value = 1541.9231
$ws.Range("J135").Value = 1252.5385
$ws.Range("K135").Value = 13877.3079
$ws.Range("L135").Value = 11272.8465
$ws.Range("M135").Value = -11342.3079
$ws.Range("N135").Value = -16342.8465

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 558.7059
$ws.Range("I2").Value = 505.57144
$ws.Range("J2").Value = 674.9375
$ws.Range("K2").Value = 505.57144
$ws.Range("L2").Value = 674.9375
$ws.Range("M2").Value = -392.57144
$ws.Range("N2").Value = -900.9375
$ws.Range("H32").Value = 12481.927
$ws.Range("I32").Value = 10738.904
$ws.Range("K32").Value = 10738.904
$ws.Range("M32").Value = -10451.904
$ws.Range("H34").Value = 16023.333
$ws.Range("I34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("M34").ClearContents()
$ws.Range("H45").Value = 1726.8334
$ws.Range("I45").Value = 1967.8334
$ws.Range("J45").Value = 1244.8334
$ws.Range("K45").Value = 1967.8334
$ws.Range("L45").Value = 1244.8334
$ws.Range("M45").Value = -1590.8334
$ws.Range("N45").Value = -1998.8334
$ws.Range("H61").Value = 7248263.5
$ws.Range("I61").Value = 9435959
$ws.Range("J61").Value = 1519.875
$ws.Range("K61").Value = 9435959
$ws.Range("L61").Value = 1519.875
$ws.Range("M61").Value = -9435747
$ws.Range("N61").Value = -1943.875
$ws.Range("H88").Value = 1556
$ws.Range("I88").Value = 1191.2
$ws.Range("J88").Value = 1860
$ws.Range("K88").Value = 1191.2
$ws.Range("L88").Value = 1860
$ws.Range("M88").Value = -785.2
$ws.Range("N88").Value = -2672
$ws.Range("H91").Value = 1556
$ws.Range("I91").Value = 1191.2
$ws.Range("J91").Value = 1860
$ws.Range("K91").Value = 1191.2
$ws.Range("L91").Value = 1860
$ws.Range("M91").Value = 212.8
$ws.Range("N91").Value = -4668
$ws.Range("H116").Value = 558.7059
$ws.Range("I116").Value = 505.57144
$ws.Range("J116").Value = 674.9375
$ws.Range("K116").Value = 505.57144
$ws.Range("L116").Value = 674.9375
$ws.Range("M116").Value = 1788.42856
$ws.Range("N116").Value = -5262.9375
$ws.Range("H122").Value = 10486.5
$ws.Range("I122").Value = 10915.4
$ws.Range("J122").Value = 9199.799999999999
$ws.Range("K122").Value = 32746.2
$ws.Range("L122").Value = 27599.4
$ws.Range("M122").Value = -30296.2
$ws.Range("N122").Value = -32499.4
$ws.Range("H136").Value = 7248263.5
$ws.Range("I136").Value = 9435959
$ws.Range("J136").Value = 1519.875
$ws.Range("K136").Value = 28307877
$ws.Range("L136").Value = 4559.625
$ws.Range("M136").Value = -28305327
$ws.Range("N136").Value = -9659.625

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 558.7059
$ws.Range("I3").Value = 505.57144
$ws.Range("J3").Value = 674.9375
$ws.Range("K3").Value = 505.57144
$ws.Range("L3").Value = 674.9375
$ws.Range("M3").Value = -391.57144
$ws.Range("N3").Value = -902.9375
$ws.Range("H99").Value = 979.2857
$ws.Range("I99").Value = 1001.25
$ws.Range("J99").Value = 950
$ws.Range("K99").Value = 1001.25
$ws.Range("L99").Value = 950
$ws.Range("M99").Value = 496.75
$ws.Range("N99").Value = -3946

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 946.6667
$ws.Range("I16").Value = 750.125
$ws.Range("J16").Value = 1171.2858
$ws.Range("K16").Value = 750.125
$ws.Range("L16").Value = 1171.2858
$ws.Range("M16").Value = -463.125
$ws.Range("N16").Value = -1745.2858
$ws.Range("H94").Value = 3548.8262
$ws.Range("I94").Value = 2956
$ws.Range("J94").Value = 3808.1875
$ws.Range("K94").Value = 2956
$ws.Range("L94").Value = 3808.1875
$ws.Range("M94").Value = -2505
$ws.Range("N94").Value = -4710.1875
$ws.Range("H99").Value = 1530.3
$ws.Range("I99").Value = 1337.625
$ws.Range("J99").Value = 2301
$ws.Range("K99").Value = 1337.625
$ws.Range("L99").Value = 2301
$ws.Range("M99").Value = 160.375
$ws.Range("N99").Value = -5297
$ws.Range("H113").Value = 946.6667
$ws.Range("I113").Value = 750.125
$ws.Range("J113").Value = 1171.2858
$ws.Range("K113").Value = 750.125
$ws.Range("L113").Value = 1171.2858
$ws.Range("M113").Value = 1419.875
$ws.Range("N113").Value = -5511.2858
$ws.Range("H126").Value = 1530.3
$ws.Range("I126").Value = 1337.625
$ws.Range("J126").Value = 2301
$ws.Range("K126").Value = 4012.875
$ws.Range("L126").Value = 6903
$ws.Range("M126").Value = -1542.875
$ws.Range("N126").Value = -11843
$ws.Range("H132").Value = 10206058
$ws.Range("I132").Value = 14707569
$ws.Range("J132").Value = 2634.9333
$ws.Range("K132").Value = 44122707
$ws.Range("L132").Value = 7904.7999
$ws.Range("M132").Value = -44120177
$ws.Range("N132").Value = -12964.7999

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H40").Value = 9230.615
$ws.Range("J40").Value = 9230.615
$ws.Range("L40").Value = 9230.615
$ws.Range("N40").Value = -9532.615
$ws.Range("H69").Value = 31485.857
$ws.Range("I69").Value = 10000
$ws.Range("J69").Value = 35066.832
$ws.Range("K69").Value = 10000
$ws.Range("L69").Value = 35066.832
$ws.Range("M69").Value = -9251
$ws.Range("N69").Value = -36564.832
$ws.Range("H72").Value = 31485.857
$ws.Range("I72").Value = 10000
$ws.Range("J72").Value = 35066.832
$ws.Range("K72").Value = 30000
$ws.Range("L72").Value = 105200.496
$ws.Range("M72").Value = -26256
$ws.Range("N72").Value = -112688.496
$ws.Range("H126").Value = 2758.6667
$ws.Range("I126").Value = 1372.4
$ws.Range("J126").Value = 3451.8
$ws.Range("K126").Value = 4117.200000000001
$ws.Range("L126").Value = 10355.4
$ws.Range("M126").Value = -1647.200000000001
$ws.Range("N126").Value = -15295.4

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 205.14815
$ws.Range("I55").Value = 63.8
$ws.Range("J55").Value = 237.27272
$ws.Range("K55").Value = 63.8
$ws.Range("L55").Value = 237.27272
$ws.Range("M55").Value = 109.2
$ws.Range("N55").Value = -583.2727199999999
$ws.Range("H136").Value = 10207766
$ws.Range("I136").Value = 11365099
$ws.Range("K136").Value = 34095297
$ws.Range("M136").Value = -34092747

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 16000
$ws.Range("J40").Value = 16000
$ws.Range("L40").Value = 16000
$ws.Range("N40").Value = -16298
$ws.Range("H132").Value = 2153.5
$ws.Range("I132").Value = 1838.6666
$ws.Range("J132").Value = 2783.1667
$ws.Range("K132").Value = 5515.9998
$ws.Range("L132").Value = 8349.500100000001
$ws.Range("M132").Value = -2985.9998
$ws.Range("N132").Value = -13409.5001

